$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test-data values in row 2 (cnic / mobile_no columns)
$ws.Range("B2").Value = "3640211897773"
$ws.Range("C2").Value = "03006943677"

# Widen columns B (cnic) and C (mobile_no) to fit the new values
$ws.Columns.Item(2).ColumnWidth = 13.25
$ws.Columns.Item(3).ColumnWidth = 11.17

# Move the active selection to D14 (as last used by the editor)
$ws.Range("D14").Select()
